# Updates cryptos list values per the Sun Jul 14 10:31:35 UTC 2024 GitHub Actions refresh.
# Values that look like plain numbers (e.g. "538.76") are written with a leading
# apostrophe so Excel keeps them stored as text (matching the source data, which is
# all text/inlineStr), then the cell style is reset to "Normal" so no stray
# quote-prefix formatting is left on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.151.03'
$ws.Range("E2").Value = '  +2.49%  '
$ws.Range("D3").Value = '3.208.69'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''538.76'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = '''146.18'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.16%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '''0.531'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.95%  '
$ws.Range("E9").Value = '  +0.29%  '
$ws.Range("E10").Value = '  +3.54%  '
$ws.Range("D11").Value = '''0.434'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.66%  '
$ws.Range("D12").Value = '3.766.66'
$ws.Range("E12").Value = '  +1.56%  '
$ws.Range("E13").Value = '  -1.18%  '
$ws.Range("D14").Value = '''26.24'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.06%  '
$ws.Range("E15").Value = '  +2.62%  '
$ws.Range("D16").Value = '60.228.42'
$ws.Range("E16").Value = '  +2.54%  '
$ws.Range("D17").Value = '3.223.07'
$ws.Range("E17").Value = '  +2.06%  '
$ws.Range("E18").Value = '  +0.74%  '
$ws.Range("D19").Value = '''13.22'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.09%  '
$ws.Range("D20").Value = '''8.36'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.02%  '
$ws.Range("D21").Value = '''380.04'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.24%  '
$ws.Range("E22").Value = '  +0.34%  '
$ws.Range("D23").Value = '''0.526'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.56%  '
$ws.Range("D24").Value = '''70.19'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.38%  '
$ws.Range("D25").Value = '''8.97'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +10.29%  '
$ws.Range("E26").Value = '  +1.23%  '
$ws.Range("D27").Value = '''1.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("D28").Value = '0.0₃0907'
$ws.Range("E28").Value = '  +3.24%  '
$ws.Range("E29").Value = '  +0.49%  '
$ws.Range("E30").Value = '  +0.45%  '
$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").Value = '''22.39'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.94%  '
$ws.Range("B32").Value = 'NEARProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D32").Value = '''5.44'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.80%  '
$ws.Range("E33").Value = '  +3.50%  '
$ws.Range("E34").Value = '  +6.06%  '
$ws.Range("D35").Value = '''156.65'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.87%  '
$ws.Range("E36").Value = '  -1.04%  '
$ws.Range("D37").Value = '2.806.92'
$ws.Range("E37").Value = '  +6.11%  '
$ws.Range("D38").Value = '''25.72'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.71%  '
$ws.Range("D39").Value = '''0.0706'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.39%  '
$ws.Range("E40").Value = '  +0.26%  '
$ws.Range("E41").Value = '  +0.46%  '
$ws.Range("D42").Value = '''39.79'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.82%  '
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").Value = '''0.0293'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.91%  '
$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D44").Value = '''0.719'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.51%  '
$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").Value = '''0.105'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.88%  '
$ws.Range("B46").Value = 'RenzoRestakedETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D46").Value = '3.252.17'
$ws.Range("E46").Value = '  +1.40%  '
$ws.Range("E47").Value = '  +1.92%  '
$ws.Range("E48").Value = '  -0.62%  '
$ws.Range("D49").Value = '''0.810'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +6.72%  '
$ws.Range("D50").Value = '''20.79'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.35%  '
$ws.Range("E51").Value = '  +0.00%  '
